$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 to the new CIFAR test accuracy value
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 10.65

# Delete rows 3 through 11 (previously rows for communication rounds 2-10)
$ws.Range("A3:B11").EntireRow.Delete()
